# Add a new list item after the last bullet in the "Questions" list:
# "Do we need another layer of security that is performs authentication?
#  (see the likes of Spotify and their authentication set up)"
#
# The new paragraph must keep the same ListParagraph style / numPr (ilvl 0,
# numId 2) and en-US language as its neighbours, and the sentence is split
# across four runs exactly like the target markup.

$d = $word.ActiveDocument

# Locate the last paragraph of the document (the "What happens to
# employees' passwords..." bullet) and add a new paragraph right after it.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.InsertParagraphAfter()

# The freshly-created (empty) paragraph inherited the ListParagraph style
# and numPr automatically; replace its contents with the exact OOXML for
# the new bullet so the run-split and <w:lastRenderedPageBreak/> marker
# match the target precisely.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Do we need another layer of security that is performs authentication?</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> (see the likes of </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Spotify</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> and their authentication set up)</w:t>
  </w:r>
</w:p>
'@

[void]$newRange.InsertXML($newParaXml)
